$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to snake_case machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix inconsistent capitalization of connector words (de/del/y/el/la/los)
# within Mexican state/municipality names so they are consistently
# Title Cased (e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga")
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B13").Value = "Amatenango De La Frontera"
$ws.Range("B29").Value = "Guadalupe Y Calvo"
$ws.Range("B31").Value = "Hidalgo Del Parral"
$ws.Range("A54").Value = "Ciudad De México"
$ws.Range("B68").Value = "Coneto De Comonfort"
$ws.Range("B75").Value = "Nombre De Dios"
$ws.Range("B80").Value = "San Juan Del Río"
$ws.Range("A84").Value = "Estado De México"
$ws.Range("B84").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B88").Value = "Ecatepec De Morelos"
$ws.Range("B95").Value = "Tlalnepantla De Baz"
$ws.Range("B103").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B109").Value = "San Francisco Del Rincón"
$ws.Range("B111").Value = "San Luis De La Paz"
$ws.Range("B112").Value = "Silao De La Victoria"
$ws.Range("B113").Value = "Valle De Santiago"
$ws.Range("B116").Value = "Acapulco De Juárez"
$ws.Range("B120").Value = "Coyuca De Benítez"
$ws.Range("B124").Value = "Técpan De Galeana"
$ws.Range("B127").Value = "Atotonilco El Grande"
$ws.Range("B134").Value = "Autlán De Navarro"
$ws.Range("B140").Value = "Huejuquilla El Alto"
$ws.Range("B141").Value = "Lagos De Moreno"
$ws.Range("B146").Value = "San Diego De Alejandría"
$ws.Range("B147").Value = "San Juan De Los Lagos"
$ws.Range("B149").Value = "San Miguel El Alto"
$ws.Range("B151").Value = "Tamazula De Gordiano"
$ws.Range("B152").Value = "Tepatitlán De Morelos"
$ws.Range("B153").Value = "Tlajomulco De Zúñiga"
$ws.Range("B156").Value = "Unión De Tula"
$ws.Range("B158").Value = "Zacoalco De Torres"
$ws.Range("B159").Value = "Zapotlán El Grande"
$ws.Range("B187").Value = "Santa María Del Oro"
$ws.Range("B194").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B196").Value = "Oaxaca De Juárez"
$ws.Range("B197").Value = "Ocotlán De Morelos"
$ws.Range("B208").Value = "Tlacolula De Matamoros"
$ws.Range("B211").Value = "Chalchicomula De Sesma"
$ws.Range("B221").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B228").Value = "Cadereyta De Montes"
$ws.Range("B233").Value = "San Ciro De Acosta"
$ws.Range("B235").Value = "Santa María Del Río"
$ws.Range("B236").Value = "Villa De Ramos"
$ws.Range("B267").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B268").Value = "Muñoz De Domingo Arenas"
$ws.Range("B272").Value = "Tepetitla De Lardizábal"
$ws.Range("B286").Value = "Martínez De La Torre"
$ws.Range("B303").Value = "Nochistlán De Mejía"
$ws.Range("B306").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B308").Value = "Villa De Cos"

# Remove trailing metadata/footer rows (313-317), which also shrinks the
# used range/dimension from A1:D317 down to A1:D311
$ws.Range("A313:A317").EntireRow.Delete()
